$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Re-append the complete grammar table a second time (rows 75-90),
# interleaving the two rule updates (rows 59 and 62) at the point where
# their revised wording was introduced so that shared-string allocation
# order matches the authored workbook.

$ws.Range("F75").Value = "<INICIO>"
$ws.Range("G75").Value = "::="
$ws.Range("H75").Value = "S_LA <CONTENIDO> "
$ws.Range("F75:H75").Style = "Bueno"

$ws.Range("F76").Value = "<CONTENIDO>"
$ws.Range("G76").Value = "::="
$ws.Range("H76").Value = "<CONJUNTO>  <CONTENIDOR> | <ER> <CONTENIDOR>"
$ws.Range("F76:H76").Style = "Bueno"

$ws.Range("F77").Value = "<CONJUNTO>"
$ws.Range("G77").Value = "::="
$ws.Range("H77").Value = "CONJ  S_DOTS  IDENTIFICADOR  S_ARROW  <DEFCONJ>  "
$ws.Range("F77:H77").Style = "Bueno"

$ws.Range("F78").Value = "<DEFCONJ>"
$ws.Range("G78").Value = "::="
$ws.Range("H78").Value = "<SEPCOMAS> | <RANGO>"
$ws.Range("F78:H78").Style = "Bueno"

$ws.Range("F79").Value = "<SEPCOMAS>"
$ws.Range("G79").Value = "::="
$ws.Range("H79").Value = "<DATOSEP>  <SEPCOMASR>"
$ws.Range("F79:H79").Style = "Bueno"

$ws.Range("F80").Value = "<DATOSEP>"
$ws.Range("G80").Value = "::="
$ws.Range("H80").Value = "PHRASE  |  NUMBER  |  LETTER"
$ws.Range("F80:H80").Style = "Bueno"

$ws.Range("F81").Value = "<SEPCOMASR>"
$ws.Range("G81").Value = "::="
$ws.Range("H81").Value = "S_COLON  <SEPCOMAS> | S_SEMICOLON"
$ws.Range("F81:H81").Style = "Bueno"

$ws.Range("F82").Value = "<RANGO>"
$ws.Range("G82").Value = "::="
$ws.Range("H82").Value = "FULLRANK S_SEMICOLON"
$ws.Range("F82:H82").Style = "Bueno"

# Update the two rules whose wording changed
$ws.Range("H62").Value = " NUMBER | LETTER | RANGE | SPACE | S_ASTERISK | S_PLUS | S_COLON | S_DOT |S_DOTS | S_SEMICOLON | S_QMARK | S_LA | S_LINE | S_LLC |S_DQUOTES | S_QUOTE"
$ws.Range("H59").Value = " NUMBER  |  LETTER | RANGE | S_ASTERISK | S_PLUS | S_COLON | S_DOT  |  S_DOTS  | S_SEMICOLON |  S_QMARK  | S_LA  |  S_LINE  |  S_LLC |S_DQUOTES | S_QUOTE S_LA  |  S_LINE  |  S_LLC"

$ws.Range("F83").Value = "<DATORANGO>"
$ws.Range("G83").Value = "::="
$ws.Range("H83").Value = "RANGE | SPACE | S_ASTERISK | S_PLUS | S_COLON | S_DOT  |  S_DOTS  | S_SEMICOLON |  S_QMARK  | S_LA  |  S_LINE  |  S_LLC"
$ws.Range("F83:H83").Style = "Bueno"

$ws.Range("F84").Value = "<ER>"
$ws.Range("G84").Value = "::="
$ws.Range("H84").Value = "IDENTIFICADOR  S_ARROW  <DEFER>   S_SEMICOLON"
$ws.Range("F84:H84").Style = "Bueno"

$ws.Range("F85").Value = "<DEFER> "
$ws.Range("G85").Value = "::="
$ws.Range("H85").Value = "<OP> | <REFCONJ> |  PHRASE | S_LBREAK |  S_QUOTE |  S_DQUOTES"
$ws.Range("F85:H85").Style = "Bueno"

$ws.Range("F86").Value = "<OP>"
$ws.Range("G86").Value = "::="
$ws.Range("H86").Value = "S_LINE <DEFER>  | S_ASTERISK  <DEFER>  | S_PLUS <DEFER>  |  S_QMARK <DEFER>  | S_DOT <DEFER> "
$ws.Range("F86:H86").Style = "Bueno"

$ws.Range("F87").Value = "<REFCONJ>"
$ws.Range("G87").Value = "::="
$ws.Range("H87").Value = "S_LA   IDENTIFICADOR   S_LLC"
$ws.Range("F87:H87").Style = "Bueno"

$ws.Range("F88").Value = "<CONTENIDOR>"
$ws.Range("G88").Value = "::="
$ws.Range("H88").Value = "<CONTENIDO> | S_PCENTS  <CADENAS>"
$ws.Range("F88:H88").Style = "Bueno"

$ws.Range("F89").Value = "<CADENAS>"
$ws.Range("G89").Value = "::="
$ws.Range("H89").Value = "IDENTIFICADOR  S_DOTS  PHRASE S_SEMICOLON <CADENASR>"
$ws.Range("F89:H89").Style = "Bueno"

$ws.Range("F90").Value = "<CADENASR>"
$ws.Range("G90").Value = "::="
$ws.Range("H90").Value = "S_LLC | <CADENAS>"
$ws.Range("F90:H90").Style = "Bueno"

# Update selection / active cell to reflect the new extent of data
$ws.Range("H91").Select()
